$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows 21-24 (rows beyond the new data range)
$ws.Rows("21:24").Delete()

# Helper: set a cell as literal text, avoiding Excel auto-converting
# numeric-looking / date-looking strings into Number/Date types, and
# avoid leaving a stray NumberFormat style behind on the cell.
function Set-TextCell($addr, $val) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = "Normal"
}


Set-TextCell "A2" "MALBEC"
Set-TextCell "B2" "58520465"
Set-TextCell "E2" "05/17/2025"
Set-TextCell "G2" "1y 2m 8d"
$ws.Range("I2").Value = 65.5

Set-TextCell "A3" "BUCKEYE"
Set-TextCell "B3" "58514522"
Set-TextCell "D3" "Cat Adoption Room G"
Set-TextCell "E3" "05/19/2025"
Set-TextCell "F3" "Domestic Shorthair"
Set-TextCell "G3" "2y 2m 2d"
Set-TextCell "H3" "Available"
$ws.Range("I3").Value = 63.5

Set-TextCell "A4" "KATNISS"
Set-TextCell "B4" "58517968"
Set-TextCell "E4" "05/16/2025"
Set-TextCell "G4" "3m 3d"
$ws.Range("I4").Value = 66.3

Set-TextCell "A5" "PRIMROSE"
Set-TextCell "B5" "58517971"
Set-TextCell "E5" "05/16/2025"
Set-TextCell "G5" "3m 3d"
$ws.Range("I5").Value = 66.3

Set-TextCell "A6" "HAYMITCH"
Set-TextCell "B6" "58517973"
Set-TextCell "E6" "05/16/2025"
Set-TextCell "F6" "Domestic Shorthair"
Set-TextCell "G6" "3m 3d"
$ws.Range("I6").Value = 66.3

Set-TextCell "A7" "PEETA"
Set-TextCell "B7" "58517974"
Set-TextCell "D7" "Foster Home"
Set-TextCell "E7" "05/16/2025"
Set-TextCell "F7" "Domestic Shorthair"
Set-TextCell "G7" "3m 3d"
Set-TextCell "H7" "In Foster"
$ws.Range("I7").Value = 66.3

Set-TextCell "A8" "Sugar"
Set-TextCell "B8" "58433959"
Set-TextCell "D8" "Offsite Adoptions"
Set-TextCell "E8" "05/05/2025"
Set-TextCell "F8" "Domestic Shorthair"
Set-TextCell "G8" "3m 0d"
Set-TextCell "H8" "Hold - Adopted!"
$ws.Range("I8").Value = 77.6

Set-TextCell "A9" "LUNA"
Set-TextCell "B9" "58067302"
Set-TextCell "C9" "Dog"
Set-TextCell "D9" "Dog Adoptions D"
Set-TextCell "E9" "07/03/2025"
Set-TextCell "F9" "Mixed Breed, Large (over 44 lbs fully grown)"
Set-TextCell "G9" "10y 3m 24d"
Set-TextCell "H9" "Available - Doggie Entourage"
$ws.Range("I9").Value = 18.5

Set-TextCell "A10" "Glow"
Set-TextCell "B10" "58834563"
Set-TextCell "C10" "Dog"
Set-TextCell "D10" "Dog Holding E"
Set-TextCell "E10" "07/02/2025"
Set-TextCell "F10" "Bulldog"
Set-TextCell "G10" "1y 6m 18d"
Set-TextCell "H10" "Hold - Cruelty Foster"
$ws.Range("I10").Value = 19.2

Set-TextCell "A11" "Katniss"
Set-TextCell "B11" "58834486"
Set-TextCell "D11" "Dog Holding E"
Set-TextCell "E11" "07/02/2025"
Set-TextCell "F11" "Bulldog"
$ws.Range("G11").ClearContents()
Set-TextCell "H11" "Hold - Cruelty Foster"
$ws.Range("I11").Value = 19.2

Set-TextCell "A12" "Sky"
Set-TextCell "B12" "58854396"
Set-TextCell "D12" "Dog Holding F"
Set-TextCell "E12" "07/07/2025"
Set-TextCell "G12" "1y 14d"
Set-TextCell "H12" "Hold - For RTO"
$ws.Range("I12").Value = 14.5

Set-TextCell "A13" "Scruffles"
Set-TextCell "B13" "58419285"
Set-TextCell "D13" "Foster Home"
Set-TextCell "E13" "05/01/2025"
Set-TextCell "F13" "Mixed Breed, Small (under 24 lbs fully grown)"
Set-TextCell "G13" "15y 2m 19d"
Set-TextCell "H13" "In Foster"
$ws.Range("I13").Value = 81.1

Set-TextCell "A14" "Dior"
Set-TextCell "B14" "58834490"
Set-TextCell "D14" "Foster Home"
Set-TextCell "E14" "07/02/2025"
Set-TextCell "F14" "Bulldog"
Set-TextCell "G14" "1y 18d"
Set-TextCell "H14" "In Foster"
$ws.Range("I14").Value = 19.2

Set-TextCell "A15" "Remy"
Set-TextCell "B15" "58834525"
Set-TextCell "D15" "Foster Home"
Set-TextCell "E15" "07/02/2025"
Set-TextCell "F15" "Bulldog, French"
Set-TextCell "G15" "1y 17d"
Set-TextCell "H15" "In Foster"
$ws.Range("I15").Value = 19.2

Set-TextCell "A16" "Bruno"
Set-TextCell "B16" "58849570"
Set-TextCell "D16" "Foster Home"
Set-TextCell "E16" "07/05/2025"
Set-TextCell "F16" "Mixed Breed, Large (over 44 lbs fully grown)"
Set-TextCell "G16" "5m 16d"
Set-TextCell "H16" "In Foster"
$ws.Range("I16").Value = 16.3

Set-TextCell "A17" "Lilly"
Set-TextCell "B17" "58831432"
Set-TextCell "E17" "07/02/2025"
Set-TextCell "F17" "Chihuahua, Long Coat"
Set-TextCell "G17" "10y 19d"
Set-TextCell "H17" "In If the Fur Fits - Medical"
$ws.Range("I17").Value = 19.4

Set-TextCell "A18" "COLT"
Set-TextCell "B18" "58838875"
Set-TextCell "E18" "07/07/2025"
Set-TextCell "F18" "Mixed Breed, Large (over 44 lbs fully grown)"
Set-TextCell "G18" "10y 13d"
Set-TextCell "H18" "In If the Fur Fits - Medical"
$ws.Range("I18").Value = 14.5

Set-TextCell "A19" "Mabel"
Set-TextCell "B19" "58421981"
Set-TextCell "C19" "Dog"
Set-TextCell "D19" "If The Fur Fits"
Set-TextCell "E19" "07/03/2025"
Set-TextCell "F19" "Mastiff"
Set-TextCell "G19" "5m 19d"
Set-TextCell "H19" "In If the Fur Fits - Trial"
$ws.Range("I19").Value = 18.5

Set-TextCell "A20" "SPRITZ"
Set-TextCell "B20" "58831524"
Set-TextCell "C20" "Rabbit"
Set-TextCell "D20" "Foster Home"
Set-TextCell "E20" "07/02/2025"
Set-TextCell "F20" "Holland Lop"
$ws.Range("G20").ClearContents()
Set-TextCell "H20" "In Foster"
$ws.Range("I20").Value = 19.4
